$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation needs to be inserted right before the
# existing row 113 (old rows 113-115 shift down to 114-116).
$ws.Rows.Item(113).Insert()

$ws.Cells.Item(113, 1).Value = 10
$ws.Cells.Item(113, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(113, 3).Value = "La Araucanía"
$ws.Cells.Item(113, 4).Value = 45239
$ws.Cells.Item(113, 5).Value = 9
$ws.Cells.Item(113, 6).Value = 300000000
$ws.Cells.Item(113, 7).Value = "Espárragos"
$ws.Cells.Item(113, 8).Value = "Sin especificar"
$ws.Cells.Item(113, 9).Value = "Primera"
$ws.Cells.Item(113, 10).Value = 500
$ws.Cells.Item(113, 11).Value = 2200
$ws.Cells.Item(113, 12).Value = 2200
$ws.Cells.Item(113, 13).Value = 2200
$ws.Cells.Item(113, 14).Value = "`$/kilo"
$ws.Cells.Item(113, 15).Value = "Región del Maule"
$ws.Cells.Item(113, 16).Value = 2200
$ws.Cells.Item(113, 17).Value = 1
$ws.Cells.Item(113, 18).Value = "Hortaliza"
